$wb = $excel.ActiveWorkbook
$wsIn = $wb.Worksheets.Item("ProductLoanInput")
$wsOut = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "4291-MS-EI-DB-SAR-REC-RNI-INT-FFConMONTHLYonDAY25-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-1st"

# --- ProductLoanInput sheet ---
$wsIn.Activate()

# productname text tweaked (ONTIME-PER -> 1st)
$wsIn.Range("B1").Value = $newProductName

# shortname changed from the number 4291 to the text "429a"
$wsIn.Range("B2").Value = "429a"

# cursor moved off the old B19 selection, no longer the active tab
$wsIn.Range("B3").Select()

# --- ProductLoanOutput sheet ---
$wsOut.Activate()

# mirror the same productname text update
$wsOut.Range("B1").Value = $newProductName

# this sheet is now the active / selected tab
$wsOut.Range("B1").Select()
